$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    This literal appears (as one shared string) in the "Status"-like cells
#    of all three sheets: Overview!E2/F2 (per-language status) and the
#    "Status" column (C2) of the zh-cn / de-de detail sheets.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the "Status" columns to match the shorter text.
#    Target raw column width (OOXML <col width>) is 13.4101845877511;
#    ColumnWidth is expressed in characters and gets snapped to whole
#    pixels, so 12.5 characters is the closest settable value that lands on
#    the nearest achievable width.
# ---------------------------------------------------------------------------
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
